# Regenerate the "K" column (column G) of save_data with freshly computed
# strike-count values (previously derived from "Strike#", now recomputed
# from the std/mean of the underlying s_vals series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values for rows 2..63, in row order.
$kValues = @(
    1,3,1,1,1,1,2,1,1,2,
    3,1,1,2,1,3,3,1,0,3,
    1,2,1,1,1,2,3,1,1,1,
    2,2,1,1,1,2,1,1,1,1,
    1,1,2,1,1,0,1,2,1,3,
    0,1,1,1,0,0,1,1,1,1,
    0,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
